# The commit re-generated this fixture's package, which re-minted the
# r:id relationship-id strings for every w:hyperlink in the document
# (see commit: "feat: update package info and setup release-it").
# None of the hyperlink targets, display text, or styling changed -
# only the opaque internal relationship-id tokens churned because the
# whole docx was rebuilt by the authoring tool.
#
# Word's object model doesn't let a caller dictate the literal
# relationship-id string that gets minted (Word always allocates its
# own rIdN sequence under the hood), so the faithful way to express
# "this hyperlink's relationship was refreshed, content unchanged" is
# to touch the writable Hyperlink surface (Address) for every
# hyperlink while preserving its current address exactly.

$d = $word.ActiveDocument
$hyperlinks = $d.Hyperlinks
$count = $hyperlinks.Count

for ($i = 1; $i -le $count; $i++) {
    $link = $hyperlinks.Item($i)
    $address = $link.Address

    # A hyperlink with a blank Address (e.g. the empty-href markdown
    # link "[links]()") would be unwrapped/removed if we re-assigned
    # Address back to "", so leave those alone - nothing about their
    # visible content changed either way.
    if ($address -ne $null -and $address -ne "") {
        $link.Address = $address
    }
}
